$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:E51 to Text format so numeric-looking strings (e.g. "248.78", "0.656")
# are stored as text, matching the source inlineStr cells, then reset the style
# index back to the default (Normal) so only the cell VALUES change - no new
# cell-level formatting is introduced.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "37.105.93"
$ws.Range("E2").Value = "  +1.34%  "

$ws.Range("D3").Value = "2.058.99"
$ws.Range("E3").Value = "  -2.31%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "248.78"
$ws.Range("E5").Value = "  -1.63%  "

$ws.Range("D6").Value = "0.656"
$ws.Range("E6").Value = "  -0.89%  "

$ws.Range("D8").Value = "55.50"
$ws.Range("E8").Value = "  +15.38%  "

$ws.Range("D9").Value = "61.51"
$ws.Range("E9").Value = "  +3.15%  "

$ws.Range("D10").Value = "0.379"
$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("D11").Value = "0.0797"
$ws.Range("E11").Value = "  +6.74%  "

$ws.Range("E12").Value = "  +5.56%  "

$ws.Range("D13").Value = "15.15"
$ws.Range("E13").Value = "  +6.03%  "

$ws.Range("D14").Value = "2.357.86"
$ws.Range("E14").Value = "  -2.53%  "

$ws.Range("D15").Value = "0.816"
$ws.Range("E15").Value = "  -1.54%  "

$ws.Range("D16").Value = "5.24"
$ws.Range("E16").Value = "  +2.85%  "

$ws.Range("D17").Value = "2.054.68"
$ws.Range("E17").Value = "  -2.59%  "

$ws.Range("D18").Value = "37.056.08"

$ws.Range("D19").Value = "0.0₃0937"
$ws.Range("E19").Value = "  +12.46%  "

$ws.Range("D20").Value = "72.36"
$ws.Range("E20").Value = "  -1.10%  "

$ws.Range("D21").Value = "14.19"
$ws.Range("E21").Value = "  +6.90%  "

$ws.Range("D22").Value = "5.39"
$ws.Range("E22").Value = "  +4.21%  "

$ws.Range("D23").Value = "237.10"
$ws.Range("E23").Value = "  -1.49%  "

$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("E25").Value = "  -1.30%  "

$ws.Range("D26").Value = "170.53"
$ws.Range("E26").Value = "  -0.55%  "

$ws.Range("D27").Value = "9.05"
$ws.Range("E27").Value = "  -1.33%  "

$ws.Range("D28").Value = "20.19"
$ws.Range("E28").Value = "  -5.83%  "

$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").Value = "4.55"
$ws.Range("E31").Value = "  +2.11%  "

$ws.Range("E32").Value = "  +11.76%  "

$ws.Range("D33").Value = "0.0624"
$ws.Range("E33").Value = "  +3.34%  "

$ws.Range("D34").Value = "4.37"
$ws.Range("E34").Value = "  +7.46%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").Value = "0.0860"
$ws.Range("E36").Value = "  -5.09%  "

$ws.Range("D37").Value = "2.27"
$ws.Range("E37").Value = "  -3.40%  "

$ws.Range("E38").Value = "  -6.53%  "

$ws.Range("D39").Value = "1.35"
$ws.Range("E39").Value = "  +1.26%  "

$ws.Range("D40").Value = "0.104"
$ws.Range("E40").Value = "  +23.54%  "

$ws.Range("D41").Value = "18.00"
$ws.Range("E41").Value = "  +11.75%  "

$ws.Range("D42").Value = "0.0223"
$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("E43").Value = "  -3.16%  "

$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "4.41"
$ws.Range("E44").Value = "  +52.00%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "96.34"
$ws.Range("E45").Value = "  -1.75%  "

$ws.Range("E46").Value = "  +1.02%  "

$ws.Range("D47").Value = "13.58"
$ws.Range("E47").Value = "  -52.71%  "

$ws.Range("D48").Value = "2.43"
$ws.Range("E48").Value = "  +6.67%  "

$ws.Range("D49").Value = "1.296.08"
$ws.Range("E49").Value = "  -3.31%  "

$ws.Range("D50").Value = "2.92"
$ws.Range("E50").Value = "  +2.93%  "

$ws.Range("D51").Value = "6.78"
$ws.Range("E51").Value = "  -5.55%  "

$dataRange.Style = "Normal"
